$wb = $excel.ActiveWorkbook

# The sheet previously named "TestDataWellFormed_Simple" (small, 4-row sample)
# is being renamed to make room for a new, fuller "TestDataWellFormed_Simple"
# sheet (the former "TestDataWellFormed_Simple_") used for multi-device smoke
# tests.
$wsFourLine = $wb.Worksheets.Item("TestDataWellFormed_Simple")
$wsSimple   = $wb.Worksheets.Item("TestDataWellFormed_Simple_")

# Rename the small sample sheet out of the way first, then claim its old name
# for the fuller data sheet. Excel keeps the _xlnm._FilterDatabase defined
# name (which pointed at "TestDataWellFormed_Simple_") in sync automatically.
$wsFourLine.Name = "TestDataWellFormed_Simple_4Line"
$wsSimple.Name = "TestDataWellFormed_Simple"

# Move the selection on the 4-line sheet, then leave it as the inactive tab.
$wsFourLine.Activate() | Out-Null
$wsFourLine.Range("E28").Select() | Out-Null

# Make the renamed (fuller) "TestDataWellFormed_Simple" sheet the active tab.
$wsSimple.Activate() | Out-Null
